$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row number => values for columns B..O (in order)
$data = @{
    2 = @(13.61277426120842, 6.63253643051909, 4.441605768527828, 0, 29.64574716163366, 37.88143490938002, 3.716965373175691, 4.378364015722758, 12.80108417555434, 21.06256085302319, 5.991847254153829, 11.20291455416862, 7.126318790424187, 0)
    3 = @(12.76515658992918, 6.295120530274113, 4.251440499035825, 0, 29.1489133971357, 37.17186802271954, 3.944579211328018, 4.56293129091833, 12.7421816436095, 20.9299887456871, 5.911968768308787, 10.51631964442569, 6.88749178360225, 0)
    4 = @(12.21560251972122, 6.081069415679404, 4.131532847424139, 0, 28.84400963095976, 36.73534473845319, 4.089412954455761, 4.680854636802387, 12.70781911377057, 20.84956765093244, 5.861139708263009, 10.07155812714568, 6.737469940990658, 0)
    5 = @(11.98267473862701, 5.997511944390001, 4.084264832918437, 0, 28.70568044094066, 36.53289386501775, 4.150131369550905, 4.732430646574811, 12.6896766259999, 20.80669362421187, 5.839021235632571, 9.885914225793515, 6.677082068293698, 0)
    6 = @(11.94144994501103, 5.990135679401313, 4.0792300697615, 0, 28.66533437672613, 36.46931667687988, 4.160782156965286, 4.744056083456092, 12.68108211479015, 20.78693731964843, 5.834151832874274, 9.856464877722026, 6.6688899710865, 0)
    7 = @(12.20683759154816, 6.097666796060752, 4.138575332527688, 0, 28.79497201066809, 36.65135848996142, 4.0915431314992, 4.689467859397594, 12.69231209801455, 20.81455812922548, 5.857677808901886, 10.07369422920471, 6.741737906734026, 0)
    8 = @(13.31959588225362, 6.539767865406566, 4.386327402970711, 0, 29.41359982936531, 37.53235246873327, 3.796055688224437, 4.451094331024822, 12.76049136931699, 20.97174549583687, 5.960616956507103, 10.97668156110452, 7.051169310152098, 0)
    9 = @(15.27798603384086, 7.322608090387561, 4.832599753533313, 0, 30.70517765154008, 39.38712077544407, 3.253813354398937, 4.006477125975567, 12.9355248899913, 21.34798772041171, 6.153757751687486, 12.55565298714453, 7.620150926591748, 0)
    10 = @(16.56124386033349, 7.901541839970506, 5.116859321566705, 0, 31.39585535912554, 40.3577912071321, 2.906190918730782, 3.712002618497181, 12.99947665148945, 21.47312929795591, 6.245901599170113, 13.61149916542323, 7.932103311700464, 0)
    11 = @(17.02141838340334, 8.447685161781125, 5.007480520375633, 0, 29.6729919010244, 37.76943242872389, 3.595995485318507, 3.650352895977591, 12.43767163511495, 20.29604536862495, 5.996571183495417, 14.12765542674903, 7.333659185428634, 0)
    12 = @(17.15366817172373, 8.780676304353285, 4.843438301642212, 0, 28.13772624451916, 35.48606323024914, 4.753239556558378, 3.640583746792428, 11.97678040425352, 19.34924205220484, 5.853570875266666, 14.3453274183266, 6.771766787361684, 0)
    13 = @(17.0395023477858, 8.999100706428822, 4.633233224188853, 0, 26.57617302117393, 33.1637113386883, 6.056436385963122, 3.683244963635268, 11.54342187322434, 18.47543211390587, 5.782267450206935, 14.35878380152066, 6.208351697677211, 0)
    14 = @(16.84525076460859, 9.108949426434464, 4.46490611083099, 0, 25.46482732618777, 31.50598512193371, 7.012543552629137, 3.739116599129862, 11.25267573554713, 17.90074809343907, 5.777059413730163, 14.27741825302839, 5.820704872234025, 0)
    15 = @(16.74292206374331, 9.116907371109924, 4.416345220723626, 0, 25.1751848727276, 31.07074788091942, 7.240363005667216, 3.768280519330608, 11.18311605318532, 17.76507907971677, 5.779422974427339, 14.21481041912327, 5.725434381925478, 0)
    16 = @(16.23180770399823, 8.859475799329081, 4.333173847346812, 0, 25.14269099109846, 31.02916702066298, 7.084732312288218, 3.890327943550711, 11.23460500678852, 17.87172581058718, 5.748454618998203, 13.78525576170031, 5.6916184766025, 0)
    17 = @(15.93934571546352, 8.596300930491987, 4.361753757868623, 0, 25.70899746845317, 31.8810603350741, 6.435300788770856, 3.955086485944488, 11.42647655883359, 18.25101109297855, 5.719544008599121, 13.49371751359746, 5.870890982138571, 0)
    18 = @(15.81232369623712, 8.290569834618831, 4.484696969900429, 0, 26.90196107676077, 33.67439454639761, 5.340679482828643, 3.966644177803323, 11.77802619419158, 18.95479814176357, 5.735316940543232, 13.29129337498939, 6.27204552650083, 0)
    19 = @(15.82000633222868, 8.016986492754524, 4.684128872358316, 0, 28.45828039369782, 35.99196001432587, 4.125808117829249, 3.947414815591546, 12.21961059890952, 19.85276820096135, 5.850031517400506, 13.18719065634759, 6.85223147196939, 0)
    20 = @(16.21932940014378, 7.799160324691388, 5.061848229537485, 0, 31.07108146106844, 39.85995165942572, 2.999497321546817, 3.814376910183933, 12.93411976168198, 21.33261411583936, 6.211663247538382, 13.3534823963612, 7.861009507274806, 0)
    21 = @(17.19182739586649, 8.174168211791114, 5.312618962085228, 0, 31.98482270930917, 41.1792584555193, 2.697260308206179, 3.579719246562586, 13.0996623708159, 21.66904152806798, 6.343375447877405, 14.11566118903713, 8.226209854029467, 0)
    22 = @(17.81493328007288, 8.411561800326286, 5.450551286831424, 0, 32.51465450588971, 41.95273435716616, 2.515378109429379, 3.422372654071927, 13.19390131187899, 21.85962462454056, 6.411155583435768, 14.59258069903558, 8.410003986777868, 0)
    23 = @(17.49036928815684, 8.269059407891472, 5.370178591457421, 0, 32.28119153470283, 41.62403858408777, 2.610863351975073, 3.495685427685773, 13.16021928694028, 21.79556923016148, 6.378385923541952, 14.3358553944641, 8.307316693798732, 0)
    24 = @(16.21936241556473, 7.742183743433903, 5.068952419790896, 0, 31.31630982431598, 40.24142813432872, 2.981844522531913, 3.797044520888743, 13.01047148387102, 21.49656054012802, 6.24621018997242, 13.32912558882824, 7.917259565735719, 0)
    25 = @(14.76455702741205, 7.144787969423942, 4.72716570360477, 0, 30.27662945597612, 38.75030639868001, 3.398524084097906, 4.137243531920029, 12.85995668119601, 21.18542617758502, 6.097703573969254, 12.15530696195982, 7.478016590154354, 0)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # column B starts at index 2
        $ws.Cells.Item([int]$row, $col).Value = $vals[$i]
    }
}